$wb = $excel.ActiveWorkbook

# --- ALC row 4 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 747.8570999999999
$ws.Range("I4").Value = 647.2
$ws.Range("K4").Value = 647.2
$ws.Range("M4").Value = -533.2

# --- ALC row 6 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 846.1
$ws.Range("I6").Value = 658.7143
$ws.Range("J6").Value = 1283.3334
$ws.Range("K6").Value = 1976.1429
$ws.Range("L6").Value = 3850.0002
$ws.Range("M6").Value = -1864.1429
$ws.Range("N6").Value = -4074.0002

# --- ALC row 58 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1900.6666
$ws.Range("I58").Value = 1321.4445
$ws.Range("K58").Value = 3964.3335
$ws.Range("M58").Value = -3814.3335

# --- ALC row 76 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 1989.5
$ws.Range("I76").Value = 1979
$ws.Range("J76").Value = 2000
$ws.Range("K76").Value = 1979
$ws.Range("L76").Value = 2000
$ws.Range("M76").Value = -1664
$ws.Range("N76").Value = -2630

# --- ALC row 79 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 1989.5
$ws.Range("I79").Value = 1979
$ws.Range("J79").Value = 2000
$ws.Range("K79").Value = 1979
$ws.Range("L79").Value = 2000
$ws.Range("M79").Value = -887
$ws.Range("N79").Value = -4184

# --- ALC row 94 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 2634.8
$ws.Range("I94").Value = 1292
$ws.Range("K94").Value = 1292
$ws.Range("M94").Value = -841

# --- ALC row 106 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 5023.25
$ws.Range("I106").Value = 5023.25
$ws.Range("K106").Value = 5023.25
$ws.Range("M106").Value = -4392.25

# --- ALC row 125 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# --- ALC row 129 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1042.7
$ws.Range("I129").Value = 825.2222
$ws.Range("K129").Value = 2475.6666
$ws.Range("M129").Value = 2524.3334

# --- ALC row 138 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2358.5454
$ws.Range("I138").Value = 1277.8572
$ws.Range("K138").Value = 3833.5716
$ws.Range("M138").Value = 1306.4284

# --- ARM row 2 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1643.4546
$ws.Range("I2").Value = 906.3
$ws.Range("J2").Value = 9015
$ws.Range("K2").Value = 906.3
$ws.Range("L2").Value = 9015
$ws.Range("M2").Value = -793.3
$ws.Range("N2").Value = -9241

# --- ARM row 32 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 939
$ws.Range("I32").Value = 939
$ws.Range("K32").Value = 939
$ws.Range("M32").Value = -652

# --- ARM row 88 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1904.7273
$ws.Range("J88").Value = 1705.75
$ws.Range("L88").Value = 1705.75
$ws.Range("N88").Value = -2517.75

# --- ARM row 91 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1904.7273
$ws.Range("J91").Value = 1705.75
$ws.Range("L91").Value = 1705.75
$ws.Range("N91").Value = -4513.75

# --- ARM row 97 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 549.0769
$ws.Range("I97").Value = 573.1667
$ws.Range("J97").Value = 260
$ws.Range("K97").Value = 573.1667
$ws.Range("L97").Value = 260
$ws.Range("M97").Value = -77.16669999999999
$ws.Range("N97").Value = -1252

# --- ARM row 102 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2770.611
$ws.Range("I102").Value = 1290.7858
$ws.Range("J102").Value = 7950
$ws.Range("K102").Value = 1290.7858
$ws.Range("L102").Value = 7950
$ws.Range("M102").Value = 331.2141999999999
$ws.Range("N102").Value = -11194

# --- ARM row 116 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1643.4546
$ws.Range("I116").Value = 906.3
$ws.Range("J116").Value = 9015
$ws.Range("K116").Value = 906.3
$ws.Range("L116").Value = 9015
$ws.Range("M116").Value = 1387.7
$ws.Range("N116").Value = -13603

# --- ARM row 122 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1139.5385
$ws.Range("I122").Value = 1109.5
$ws.Range("K122").Value = 3328.5
$ws.Range("M122").Value = -878.5

# --- BSM row 3 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1643.4546
$ws.Range("I3").Value = 906.3
$ws.Range("J3").Value = 9015
$ws.Range("K3").Value = 906.3
$ws.Range("L3").Value = 9015
$ws.Range("M3").Value = -792.3
$ws.Range("N3").Value = -9243

# --- BSM row 94 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 278.8
$ws.Range("I94").Value = 273.625
$ws.Range("K94").Value = 273.625
$ws.Range("M94").Value = 177.375

# --- BSM row 99 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2901
$ws.Range("I99").Value = 2419.3333
$ws.Range("K99").Value = 2419.3333
$ws.Range("M99").Value = -921.3332999999998

# --- CRP row 62 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5549.6665
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# --- CRP row 65 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 5549.6665
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# --- CUL row 122 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 835.5
$ws.Range("I122").Value = 697.6667
$ws.Range("J122").Value = 918.2
$ws.Range("K122").Value = 6279.0003
$ws.Range("L122").Value = 8263.800000000001
$ws.Range("M122").Value = -3829.0003
$ws.Range("N122").Value = -13163.8

# --- CUL row 139 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 4499.75
$ws.Range("I139").Value = 2949.5
$ws.Range("K139").Value = 8848.5
$ws.Range("M139").Value = -3708.5

# --- GSM row 80 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1969
$ws.Range("I80").Value = 1756.625
$ws.Range("J80").Value = 2535.3333
$ws.Range("K80").Value = 1756.625
$ws.Range("L80").Value = 2535.3333
$ws.Range("M80").Value = -758.625
$ws.Range("N80").Value = -4531.3333

# --- GSM row 83 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 1969
$ws.Range("I83").Value = 1756.625
$ws.Range("J83").Value = 2535.3333
$ws.Range("K83").Value = 8783.125
$ws.Range("L83").Value = 12676.6665
$ws.Range("M83").Value = -3791.125
$ws.Range("N83").Value = -22660.6665

# --- GSM row 122 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2549.4614
$ws.Range("I122").Value = 2549.4614
$ws.Range("K122").Value = 7648.3842
$ws.Range("M122").Value = -5198.3842

# --- LTW row 61 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2400.3125
$ws.Range("I61").Value = 783.9167
$ws.Range("J61").Value = 7249.5
$ws.Range("K61").Value = 783.9167
$ws.Range("L61").Value = 7249.5
$ws.Range("M61").Value = -581.9167
$ws.Range("N61").Value = -7653.5

# --- LTW row 93 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 828.8333
$ws.Range("I93").Value = 794.6
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 794.6
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = 453.4
$ws.Range("N93").Value = -3496

# --- LTW row 113 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2400.3125
$ws.Range("I113").Value = 783.9167
$ws.Range("J113").Value = 7249.5
$ws.Range("K113").Value = 783.9167
$ws.Range("L113").Value = 7249.5
$ws.Range("M113").Value = 1386.0833
$ws.Range("N113").Value = -11589.5

# --- LTW row 133 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# --- WVR row 81 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 101
$ws.Range("I81").Value = 101
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 202
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 859
$ws.Range("N81").ClearContents()

# --- WVR row 84 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 101
$ws.Range("I84").Value = 101
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 1010
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 4294
$ws.Range("N84").ClearContents()

# --- WVR row 96 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1128.2
$ws.Range("I96").Value = 782.3333
$ws.Range("J96").Value = 1647
$ws.Range("K96").Value = 782.3333
$ws.Range("L96").Value = 1647
$ws.Range("M96").Value = 590.6667
$ws.Range("N96").Value = -4393

# --- WVR row 100 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 975.7857
$ws.Range("I100").Value = 753.44446
$ws.Range("K100").Value = 1506.88892
$ws.Range("M100").Value = -965.8889200000001

# --- WVR row 107 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 5198.5
$ws.Range("I107").Value = 5198.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 15595.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -13675.5
$ws.Range("N107").ClearContents()

# --- WVR row 136 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3532
$ws.Range("J136").Value = 5255.5454
$ws.Range("L136").Value = 15766.6362
$ws.Range("N136").Value = -20866.6362
